# Scheduled market-data refresh: updates currentAveragePrice* / Leve*Price / LeveProfit*
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly
# pulled prices, exactly as the runner's nightly job does.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3499.1667
$ws.Cells.Item(76, 9).Value = 3508.182
$ws.Cells.Item(76, 11).Value = 3508.182
$ws.Cells.Item(76, 13).Value = -3193.182

$ws.Cells.Item(79, 8).Value = 3499.1667
$ws.Cells.Item(79, 9).Value = 3508.182
$ws.Cells.Item(79, 11).Value = 3508.182
$ws.Cells.Item(79, 13).Value = -2416.182

$ws.Cells.Item(113, 8).Value = 3838.125
$ws.Cells.Item(113, 10).Value = 3450
$ws.Cells.Item(113, 12).Value = 3450
$ws.Cells.Item(113, 14).Value = -9958

$ws.Cells.Item(135, 8).Value = 49291.617
$ws.Cells.Item(135, 9).Value = 78367.69500000001
$ws.Cells.Item(135, 10).Value = 2043
$ws.Cells.Item(135, 11).Value = 705309.2550000001
$ws.Cells.Item(135, 12).Value = 18387
$ws.Cells.Item(135, 13).Value = -702774.2550000001
$ws.Cells.Item(135, 14).Value = -23457

$ws.Cells.Item(137, 8).Value = 2633160.5
$ws.Cells.Item(137, 9).Value = 3031738.2
$ws.Cells.Item(137, 10).Value = 2548
$ws.Cells.Item(137, 11).Value = 9095214.600000001
$ws.Cells.Item(137, 12).Value = 7644
$ws.Cells.Item(137, 13).Value = -9092664.600000001
$ws.Cells.Item(137, 14).Value = -12744

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 50101840
$ws.Cells.Item(61, 9).Value = 77001256
$ws.Cells.Item(61, 10).Value = 145782.86
$ws.Cells.Item(61, 11).Value = 77001256
$ws.Cells.Item(61, 12).Value = 145782.86
$ws.Cells.Item(61, 13).Value = -77001044
$ws.Cells.Item(61, 14).Value = -146206.86

$ws.Cells.Item(74, 8).Value = 8131038
$ws.Cells.Item(74, 9).Value = 11410550
$ws.Cells.Item(74, 10).Value = 114451.555
$ws.Cells.Item(74, 11).Value = 11410550
$ws.Cells.Item(74, 12).Value = 114451.555
$ws.Cells.Item(74, 13).Value = -11409676
$ws.Cells.Item(74, 14).Value = -116199.555

$ws.Cells.Item(77, 8).Value = 8131038
$ws.Cells.Item(77, 9).Value = 11410550
$ws.Cells.Item(77, 10).Value = 114451.555
$ws.Cells.Item(77, 11).Value = 57052750
$ws.Cells.Item(77, 12).Value = 572257.7749999999
$ws.Cells.Item(77, 13).Value = -57048382
$ws.Cells.Item(77, 14).Value = -580993.7749999999

$ws.Cells.Item(132, 8).Value = 44879.168
$ws.Cells.Item(132, 9).Value = 27772.27
$ws.Cells.Item(132, 10).Value = 102420.55
$ws.Cells.Item(132, 11).Value = 83316.81
$ws.Cells.Item(132, 12).Value = 307261.65
$ws.Cells.Item(132, 13).Value = -80786.81
$ws.Cells.Item(132, 14).Value = -312321.65

$ws.Cells.Item(136, 8).Value = 50101840
$ws.Cells.Item(136, 9).Value = 77001256
$ws.Cells.Item(136, 10).Value = 145782.86
$ws.Cells.Item(136, 11).Value = 231003768
$ws.Cells.Item(136, 12).Value = 437348.58
$ws.Cells.Item(136, 13).Value = -231001218
$ws.Cells.Item(136, 14).Value = -442448.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(141, 8).Value = 41241.6
$ws.Cells.Item(141, 10).Value = 41596.668
$ws.Cells.Item(141, 12).Value = 41596.668
$ws.Cells.Item(141, 14).Value = -51956.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 30304834
$ws.Cells.Item(58, 9).Value = 33334958
$ws.Cells.Item(58, 10).Value = 3600.3333
$ws.Cells.Item(58, 11).Value = 33334958
$ws.Cells.Item(58, 12).Value = 3600.3333
$ws.Cells.Item(58, 13).Value = -33334755
$ws.Cells.Item(58, 14).Value = -4006.3333

$ws.Cells.Item(94, 8).Value = 4357.7856
$ws.Cells.Item(94, 9).Value = 9922.200000000001
$ws.Cells.Item(94, 11).Value = 9922.200000000001
$ws.Cells.Item(94, 13).Value = -9471.200000000001

$ws.Cells.Item(105, 8).Value = 1066.3334
$ws.Cells.Item(105, 9).Value = 1100
$ws.Cells.Item(105, 10).Value = 999
$ws.Cells.Item(105, 11).Value = 1100
$ws.Cells.Item(105, 12).Value = 999
$ws.Cells.Item(105, 13).Value = 647
$ws.Cells.Item(105, 14).Value = -4493

$ws.Cells.Item(122, 8).Value = 1805.5652
$ws.Cells.Item(122, 9).Value = 1383.1765
$ws.Cells.Item(122, 11).Value = 4149.529500000001
$ws.Cells.Item(122, 13).Value = -1699.529500000001

$ws.Cells.Item(132, 8).Value = 30125.416
$ws.Cells.Item(132, 9).Value = 2192.6155
$ws.Cells.Item(132, 10).Value = 102750.7
$ws.Cells.Item(132, 11).Value = 6577.8465
$ws.Cells.Item(132, 12).Value = 308252.1
$ws.Cells.Item(132, 13).Value = -4047.8465
$ws.Cells.Item(132, 14).Value = -313312.1

$ws.Cells.Item(136, 8).Value = 30304834
$ws.Cells.Item(136, 9).Value = 33334958
$ws.Cells.Item(136, 10).Value = 3600.3333
$ws.Cells.Item(136, 11).Value = 100004874
$ws.Cells.Item(136, 12).Value = 10800.9999
$ws.Cells.Item(136, 13).Value = -100002324
$ws.Cells.Item(136, 14).Value = -15900.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 983
$ws.Cells.Item(32, 9).Value = 1650
$ws.Cells.Item(32, 10).Value = 649.5
$ws.Cells.Item(32, 11).Value = 4950
$ws.Cells.Item(32, 12).Value = 1948.5
$ws.Cells.Item(32, 13).Value = -4667
$ws.Cells.Item(32, 14).Value = -2514.5

$ws.Cells.Item(55, 8).Value = 2384.923
$ws.Cells.Item(55, 9).Value = 1004
$ws.Cells.Item(55, 10).Value = 2500
$ws.Cells.Item(55, 11).Value = 3012
$ws.Cells.Item(55, 12).Value = 7500
$ws.Cells.Item(55, 13).Value = -2835
$ws.Cells.Item(55, 14).Value = -7854

$ws.Cells.Item(92, 8).Value = 1031.25
$ws.Cells.Item(92, 9).Value = 1031.25
$ws.Cells.Item(92, 11).Value = 3093.75
$ws.Cells.Item(92, 13).Value = -1845.75

$ws.Cells.Item(121, 8).Value = 38385856
$ws.Cells.Item(121, 9).Value = 905.55554
$ws.Cells.Item(121, 10).Value = 47477028
$ws.Cells.Item(121, 11).Value = 2716.66662
$ws.Cells.Item(121, 12).Value = 142431084
$ws.Cells.Item(121, 13).Value = -1406.66662
$ws.Cells.Item(121, 14).Value = -142433704

$ws.Cells.Item(122, 8).Value = 757.6111
$ws.Cells.Item(122, 9).Value = 506.2857
$ws.Cells.Item(122, 10).Value = 917.5454999999999
$ws.Cells.Item(122, 11).Value = 4556.571300000001
$ws.Cells.Item(122, 12).Value = 8257.9095
$ws.Cells.Item(122, 13).Value = -2106.571300000001
$ws.Cells.Item(122, 14).Value = -13157.9095

$ws.Cells.Item(131, 8).Value = 996.625
$ws.Cells.Item(131, 10).Value = 1037.0137
$ws.Cells.Item(131, 12).Value = 3111.0411
$ws.Cells.Item(131, 14).Value = -13191.0411

$ws.Cells.Item(132, 8).Value = 880
$ws.Cells.Item(132, 9).Value = 642
$ws.Cells.Item(132, 10).Value = 1475
$ws.Cells.Item(132, 11).Value = 5778
$ws.Cells.Item(132, 12).Value = 13275
$ws.Cells.Item(132, 13).Value = -3248
$ws.Cells.Item(132, 14).Value = -18335

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 126718.375
$ws.Cells.Item(132, 9).Value = 91954.17999999999
$ws.Cells.Item(132, 10).Value = 203199.6
$ws.Cells.Item(132, 11).Value = 275862.54
$ws.Cells.Item(132, 12).Value = 609598.8
$ws.Cells.Item(132, 13).Value = -273332.54
$ws.Cells.Item(132, 14).Value = -614658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3514.1538
$ws.Cells.Item(122, 9).Value = 3155.4443
$ws.Cells.Item(122, 10).Value = 4321.25
$ws.Cells.Item(122, 11).Value = 9466.332900000001
$ws.Cells.Item(122, 12).Value = 12963.75
$ws.Cells.Item(122, 13).Value = -7016.332900000001
$ws.Cells.Item(122, 14).Value = -17863.75

$ws.Cells.Item(132, 8).Value = 60817.883
$ws.Cells.Item(132, 9).Value = 36184.535
$ws.Cells.Item(132, 10).Value = 94408.82000000001
$ws.Cells.Item(132, 11).Value = 108553.605
$ws.Cells.Item(132, 12).Value = 283226.46
$ws.Cells.Item(132, 13).Value = -106023.605
$ws.Cells.Item(132, 14).Value = -288286.46

$ws.Cells.Item(136, 8).Value = 129806.44
$ws.Cells.Item(136, 9).Value = 112755.89
$ws.Cells.Item(136, 11).Value = 338267.67
$ws.Cells.Item(136, 13).Value = -335717.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 34200
$ws.Cells.Item(124, 10).Value = 34200
$ws.Cells.Item(124, 12).Value = 34200
$ws.Cells.Item(124, 14).Value = -44020

$ws.Cells.Item(132, 8).Value = 66191.28999999999
$ws.Cells.Item(132, 9).Value = 60007.766
$ws.Cells.Item(132, 11).Value = 180023.298
$ws.Cells.Item(132, 13).Value = -177493.298
